$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename the header labels (also swaps which shared-string index maps to which column)
$ws.Range("A1").Value = "eenheden kolom"
$ws.Range("B1").Value = "tientallén kolom"
$ws.Range("C1").Value = "hondertallen"

# Widen column A (target raw OOXML width 28.28515625; closest reachable value
# through the ColumnWidth COM setter, which snaps to 1/6-character increments)
$ws.Columns.Item(1).ColumnWidth = 27.5

# Change the active selection on the sheet
$ws.Range("B2").Select() | Out-Null
